$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 16
$ws.Range("H16").Value = 29135.182
$ws.Range("J16").Value = 29135.182
$ws.Range("L16").Value = 29135.182
$ws.Range("N16").Value = -29595.182

# Row 17
$ws.Range("H17").Value = 2520.22
$ws.Range("J17").Value = 2520.22
$ws.Range("L17").Value = 7560.66
$ws.Range("N17").Value = -7896.66

# Row 20
$ws.Range("H20").Value = 12851.667
$ws.Range("I20").Value = 1500
$ws.Range("J20").Value = 35555
$ws.Range("K20").Value = 1500
$ws.Range("L20").Value = 35555
$ws.Range("M20").Value = -1270
$ws.Range("N20").Value = -36015

# Row 32
$ws.Range("H32").Value = 47621356
$ws.Range("J32").Value = 2691.5
$ws.Range("L32").Value = 2691.5
$ws.Range("N32").Value = -3343.5

# Row 35
$ws.Range("H35").Value = 12851.667
$ws.Range("I35").Value = 1500
$ws.Range("J35").Value = 35555
$ws.Range("K35").Value = 1500
$ws.Range("L35").Value = 35555
$ws.Range("M35").Value = -1121
$ws.Range("N35").Value = -36313

# Row 125
$ws.Range("H125").Value = 1395
$ws.Range("I125").Value = 1412.8
$ws.Range("J125").Value = 1382.2858
$ws.Range("K125").Value = 12715.2
$ws.Range("L125").Value = 12440.5722
$ws.Range("M125").Value = -10255.2
$ws.Range("N125").Value = -17360.5722

# Row 129
$ws.Range("H129").Value = 1050.6852
$ws.Range("I129").Value = 1002
$ws.Range("J129").Value = 1075.0278
$ws.Range("K129").Value = 3006
$ws.Range("L129").Value = 3225.0834
$ws.Range("M129").Value = 1994
$ws.Range("N129").Value = -13225.0834

# Row 135
$ws.Range("H135").Value = 50002748
$ws.Range("I135").Value = 2183
$ws.Range("J135").Value = 250005000
$ws.Range("K135").Value = 19647
$ws.Range("L135").Value = 2250045000
$ws.Range("M135").Value = -17112
$ws.Range("N135").Value = -2250050070

# Row 137
$ws.Range("H137").Value = 2854084.2
$ws.Range("I137").Value = 8547986
$ws.Range("K137").Value = 25643958
$ws.Range("M137").Value = -25641408

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 11628.654
$ws.Range("I32").Value = 10261.24
$ws.Range("J32").Value = 25302.8
$ws.Range("K32").Value = 10261.24
$ws.Range("L32").Value = 25302.8
$ws.Range("M32").Value = -9974.24
$ws.Range("N32").Value = -25876.8

# Row 34
$ws.Range("H34").Value = 14750
$ws.Range("I34").Value = 14500
$ws.Range("J34").Value = 15000
$ws.Range("K34").Value = 14500
$ws.Range("L34").Value = 15000
$ws.Range("M34").Value = -14229
$ws.Range("N34").Value = -15542

# Row 35
$ws.Range("H35").Value = 5000
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# Row 64
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496

# Row 67
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716

# Row 74
$ws.Range("H74").Value = 2176.4243
$ws.Range("I74").Value = 1965.8636
$ws.Range("K74").Value = 1965.8636
$ws.Range("M74").Value = -1091.8636

# Row 77
$ws.Range("H77").Value = 2176.4243
$ws.Range("I77").Value = 1965.8636
$ws.Range("K77").Value = 9829.317999999999
$ws.Range("M77").Value = -5461.317999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 8136525
$ws.Range("I31").Value = 2565.0527
$ws.Range("J31").Value = 15161309
$ws.Range("K31").Value = 2565.0527
$ws.Range("L31").Value = 15161309
$ws.Range("M31").Value = -2270.0527
$ws.Range("N31").Value = -15161899

# Row 34
$ws.Range("H34").Value = 8136525
$ws.Range("I34").Value = 2565.0527
$ws.Range("J34").Value = 15161309
$ws.Range("K34").Value = 2565.0527
$ws.Range("L34").Value = 15161309
$ws.Range("M34").Value = -2363.0527
$ws.Range("N34").Value = -15161713

# Row 39
$ws.Range("H39").Value = 16833.166
$ws.Range("I39").Value = 5249.75
$ws.Range("J39").Value = 40000
$ws.Range("K39").Value = 5249.75
$ws.Range("L39").Value = 40000
$ws.Range("M39").Value = -4858.75
$ws.Range("N39").Value = -40782

# Row 49
$ws.Range("H49").Value = 16833.166
$ws.Range("I49").Value = 5249.75
$ws.Range("J49").Value = 40000
$ws.Range("K49").Value = 5249.75
$ws.Range("L49").Value = 40000
$ws.Range("M49").Value = -5067.75
$ws.Range("N49").Value = -40364

# Row 132
$ws.Range("H132").Value = 638618
$ws.Range("I132").Value = 1452.8462
$ws.Range("J132").Value = 1558967.6
$ws.Range("K132").Value = 4358.5386
$ws.Range("L132").Value = 4676902.800000001
$ws.Range("M132").Value = -1828.5386
$ws.Range("N132").Value = -4681962.800000001

$ws = $wb.Worksheets.Item("CUL")
# Row 49
$ws.Range("H49").Value = 1500
$ws.Range("J49").Value = 1500
$ws.Range("L49").Value = 4500
$ws.Range("N49").Value = -4812

# Row 129
$ws.Range("H129").Value = 201526.14
$ws.Range("I129").Value = 750982.5
$ws.Range("J129").Value = 1723.8182
$ws.Range("K129").Value = 2252947.5
$ws.Range("L129").Value = 5171.4546
$ws.Range("M129").Value = -2247947.5
$ws.Range("N129").Value = -15171.4546

# Row 131
$ws.Range("H131").Value = 2544.0532
$ws.Range("I131").Value = 12928.625
$ws.Range("J131").Value = 1304.1045
$ws.Range("K131").Value = 38785.875
$ws.Range("L131").Value = 3912.3135
$ws.Range("M131").Value = -33745.875
$ws.Range("N131").Value = -13992.3135

# Row 132
$ws.Range("H132").Value = 3312.9614
$ws.Range("J132").Value = 4279
$ws.Range("L132").Value = 38511
$ws.Range("N132").Value = -43571

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5005.7646
$ws.Range("I70").Value = 5064.24
$ws.Range("K70").Value = 5064.24
$ws.Range("M70").Value = -4794.24

# Row 73
$ws.Range("H73").Value = 5005.7646
$ws.Range("I73").Value = 5064.24
$ws.Range("K73").Value = 5064.24
$ws.Range("M73").Value = -4128.24

# Row 139
$ws.Range("H139").Value = 29980
$ws.Range("J139").Value = 29980
$ws.Range("L139").Value = 29980
$ws.Range("N139").Value = -40260

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2763.9333
$ws.Range("I7").Value = 1820.5
$ws.Range("J7").Value = 3842.1428
$ws.Range("K7").Value = 1820.5
$ws.Range("L7").Value = 3842.1428
$ws.Range("M7").Value = -1708.5
$ws.Range("N7").Value = -4066.1428

# Row 126
$ws.Range("H126").Value = 2763.9333
$ws.Range("I126").Value = 1820.5
$ws.Range("J126").Value = 3842.1428
$ws.Range("K126").Value = 5461.5
$ws.Range("L126").Value = 11526.4284
$ws.Range("M126").Value = -2991.5
$ws.Range("N126").Value = -16466.4284

$ws = $wb.Worksheets.Item("WVR")
# Row 42
$ws.Range("H42").Value = 28000
$ws.Range("J42").Value = 28000
$ws.Range("L42").Value = 28000
$ws.Range("N42").Value = -28756

# Row 81
$ws.Range("H81").Value = 1500

# Row 84
$ws.Range("H84").Value = 1500

# Row 100
$ws.Range("H100").Value = 1025.375
$ws.Range("I100").Value = 1150
$ws.Range("J100").Value = 651.5
$ws.Range("K100").Value = 2300
$ws.Range("L100").Value = 1303
$ws.Range("M100").Value = -1759
$ws.Range("N100").Value = -2385

# Row 132
$ws.Range("H132").Value = 1244320.2
$ws.Range("I132").Value = 1813113.2
$ws.Range("J132").Value = 3317.2727
$ws.Range("K132").Value = 5439339.6
$ws.Range("L132").Value = 9951.8181
$ws.Range("M132").Value = -5436809.6
$ws.Range("N132").Value = -15011.8181
